$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section: Source Type: Statistical Institution (Most Widely Used)
# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B10").Value = "'20.11"
$ws.Range("C10").Value = "'1.37"
$ws.Range("D10").Value = "'21.49"

# Employment (% of total): SMEs / MSMEs
$ws.Range("C12").Value = "'21.76"
$ws.Range("D12").Value = "'84.36"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B14").Value = "'93.46"
$ws.Range("C14").Value = "'6.39"
$ws.Range("D14").Value = "'99.84"

# Section: Source Type: SME Associations
# Enterprises density (per 1000 people): MSMEs
$ws.Range("D29").Value = "'19.37"
